$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update data rows 2-5 with new sensor readings ---
$rowData = @{}
$rowData[2] = @(45126.50694444445,4.347,5.032,0.777,10.022,8.942,4.483,9.151999999999999,5.852,2.498,6.03,7.398,4.257,0.73,3.371,5.736,3.045,0.138,0.227,54.637,11.069,2.959,6.58,4.913,0.755,5.641,3.12,5.177,2.517,6.419,0.192,7.159,2.015,4.399)
$rowData[3] = @(45126.51388888889,17.365,13.774,0.909,37.953,31.699,14.067,48.215,21.357,9.587999999999999,15.208,16.866,16.125,4.155,13.471,19.872,11.43,0.142,0.527,204.02,38.724,12.324,25.879,14.318,1.888,25.098,11.171,10.996,11.312,17.036,0.115,43.563,7.391,15.815)
$rowData[4] = @(45126.52083333334,17.988,13.974,0.824,39.246,32.653,14.401,55.419,22.024,9.906000000000001,15.321,16.829,16.664,4.378,13.982,20.442,11.793,0.139,0.516,209.899,39.907,12.817,26.82,14.559,1.904,27.529,11.521,10.902,11.777,17.223,0.092,50.302,7.612,16.312)
$rowData[5] = @(45126.52777777778,6.52,5.24,0.39,14.28,11.96,5.31,24.85,8.1,3.68,5.82,6.54,6.06,1.54,5.03,7.58,4.34,0.04,0.17,72.64,14.87,4.57,9.949999999999999,5.54,0.71,11.94,4.19,4.28,4.25,6.52,0.08,22.75,2.86,5.94)

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $vals[$c]
    }
}

# --- Remove obsolete row 6 (trims dataset to 4 data rows) ---
$ws.Rows(6).Delete()

# --- Widen specific columns from 7 to 8 character units ---
$colsToWiden = @("C","G","K","L","M","O","Q","V","X","AA","AB","AC","AD","AH")
foreach ($col in $colsToWiden) {
    $ws.Columns($col).ColumnWidth = 7.17
}
